# Re-pull data, push all data, mean calculation
# Update column F (dSF) values for the rows whose underlying source data changed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    5  = -5
    9  = -2
    12 = -4
    14 = -5
    15 = 0
    17 = 1
    18 = -2
    21 = -1
    25 = -6
    26 = -4
    31 = -2
    32 = 1
    35 = 0
    41 = 3
    42 = -1
    43 = 2
    44 = 5
    47 = 4
    48 = 4
    49 = -4
    50 = -2
    51 = -4
    53 = 0
    54 = -4
    55 = -5
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
